$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table5")
$ws.Activate()

$ws.Range("A5").Value = "prince"
$ws.Range("B5").Value = "anto"
$ws.Range("C5").Value = "manager"

$ws.Range("F9").Select()
